# Update "想去人数" (interested-count) figures in the "展览" (sheet1),
# "本地生活" (sheet3) and "全部类型" (sheet4) worksheets to reflect the
# refreshed scrape captured in commit 456a3b4.

$wb = $excel.ActiveWorkbook

$sheetExhibit = $wb.Worksheets.Item("展览")
$sheetLocal   = $wb.Worksheets.Item("本地生活")
$sheetAll     = $wb.Worksheets.Item("全部类型")

# 展览 sheet (sheet1.xml)
$sheetExhibit.Range("F4").Value  = 384
$sheetExhibit.Range("F10").Value = 2608
$sheetExhibit.Range("F13").Value = 2624
$sheetExhibit.Range("F16").Value = 2051
$sheetExhibit.Range("F26").Value = 1270
$sheetExhibit.Range("F31").Value = 1647
$sheetExhibit.Range("F34").Value = 976

# 本地生活 sheet (sheet3.xml)
$sheetLocal.Range("F3").Value = 942

# 全部类型 sheet (sheet4.xml)
$sheetAll.Range("F4").Value  = 942
$sheetAll.Range("F7").Value  = 384
$sheetAll.Range("F15").Value = 2608
$sheetAll.Range("F18").Value = 2624
$sheetAll.Range("F21").Value = 2051
$sheetAll.Range("F31").Value = 1270
$sheetAll.Range("F36").Value = 1647
$sheetAll.Range("F39").Value = 976
